# Update cucumber scenario step text in the ScenarioName column (C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "User verify login with valid credentials"
$ws.Range("C4").Value = "Validate the functionality for creating a new client"
$ws.Range("C5").Value = "Validate the functionality for updating a client"
$ws.Range("C6").Value = "Validate the functionality for deleting a client"

# Move the active selection to F11 (matches the saved sheet view state)
$ws.Range("F11").Select()
